$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.744.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.010.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.62%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.97"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +9.52%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.996.61"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.11%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.32"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +13.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.08"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.123"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.507.09"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.21"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.004.69"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.724.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "440.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.65"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.64%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.55"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.11"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.65"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.72%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +12.62%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.82"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.91"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0790"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +17.73%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +8.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.11"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.70%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.25%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +12.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.60"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "402.95"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0355"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.771.58"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.107"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.253"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.02%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.38"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.06"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.22%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.55"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +22.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.111"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.80"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.35%  "
